$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns A:F (No, Product Name, Specification,
# Origin HSCode(China), FOB/CFR/CIF, image). We need to insert three new
# columns (Qty, Unit, Unit Price) between "Origin HSCode(China)" (D) and
# "FOB/CFR/CIF" (old E) so the final layout is:
#   A=No, B=Product Name, C=Specification, D=Origin HSCode(China),
#   E=Qty, F=Unit, G=Unit Price, H=FOB/CFR/CIF, I=image

# Insert 3 blank columns at E:G - this shifts the old E (FOB/CFR/CIF) and
# old F (image) columns to H and I, carrying their values/styles with them.
$ws.Range("E1:G1").EntireColumn.Insert()

# The newly inserted header cells (E1:G1) inherit the format of the column
# to their left (D1). Re-stripe them to match the "Qty/Unit/Unit Price"
# header look (same look as the image header, column I).
$ws.Range("I1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row labels for the new columns.
$ws.Range("E1").Value = "Qty"
$ws.Range("F1").Value = "Unit"
$ws.Range("G1").Value = "Unit Price"

# Sample data row values for the new columns.
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "pcs"
$ws.Range("G2").Value = 0

# Widen the new Qty/Unit/Unit Price columns (and the neighbouring
# Origin HSCode column D, which now visually groups with them) and
# nudge the FOB/CFR/CIF column (H) to the post-edit layout widths.
$ws.Columns("D:G").ColumnWidth = 22.916666666278616
$ws.Columns("H").ColumnWidth = 20.25000000023283

# Match the new active selection left at F1 (first unit column).
$ws.Range("F1").Select()
